$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "purpose" column (E2:E19) from "S.GISH" to "fullRNASEQ"
$ws.Range("E2:E19").Value = "fullRNASEQ"

# Move selection to match the author's final cursor position after the edit
$ws.Range("D20:F23").Select()
